# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 13
$wsExhibit.Range("F4").Value = 168
$wsExhibit.Range("F5").Value = 2783
$wsExhibit.Range("F9").Value = 121
$wsExhibit.Range("F10").Value = 63
$wsExhibit.Range("F11").Value = 72
$wsExhibit.Range("F12").Value = 2601
$wsExhibit.Range("F13").Value = 765

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13
$wsAll.Range("F5").Value = 168
$wsAll.Range("F6").Value = 2783
$wsAll.Range("F11").Value = 121
$wsAll.Range("F12").Value = 63
$wsAll.Range("F13").Value = 72
$wsAll.Range("F14").Value = 2601
$wsAll.Range("F15").Value = 765
